# Apply the HIPAA employee status ValueSet metadata update:
#  - Version bump 5.0.0 -> 6.0.0
#  - Date refresh
#  - Publisher value filled in ("Alvearie Team")
#  - "Contact" / "No display for ContactDetail" duplicate rows collapsed into
#    a single "Jurisdiction" / "United States of America" row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version
$ws.Range("B3").Value = "6.0.0"

# Date
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was blank)
$ws.Range("B9").Value = "Alvearie Team"

# Replace the first "Contact" row with "Jurisdiction"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Remove the duplicate "Contact" / "No display for ContactDetail" row
$ws.Rows.Item(11).Delete()
